# Auto-generated Excel COM-interop script
# Applies a scheduled market-data refresh to the per-job profit tables
# across all eight job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
#
# For each changed row, columns H/I/J/K/L/M/N hold the refreshed
# currentAveragePrice* / LevePrice* / LeveProfit* figures. A few cells
# that previously held a value now come back empty from the price feed
# (no HQ profit data available) and are cleared instead of zeroed; one
# cell (ARM!N3) gains a value it didn't have before.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4821.5884
$ws.Range("J64").Value = 5000
$ws.Range("L64").Value = 5000
$ws.Range("N64").Value = -5496
$ws.Range("H67").Value = 4821.5884
$ws.Range("J67").Value = 5000
$ws.Range("L67").Value = 5000
$ws.Range("N67").Value = -6716
$ws.Range("H70").Value = 2000
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 2000
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H111").Value = 2275.5715
$ws.Range("I111").Value = 976.3333
$ws.Range("K111").Value = 2928.9999
$ws.Range("M111").Value = 138.0001000000002
$ws.Range("H116").Value = 5912.2915
$ws.Range("I116").Value = 5052.8237
$ws.Range("K116").Value = 5052.8237
$ws.Range("M116").Value = -1610.8237
$ws.Range("H137").Value = 261692.52
$ws.Range("I137").Value = 591329.4
$ws.Range("K137").Value = 1773988.2
$ws.Range("M137").Value = -1771438.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1370.4546
$ws.Range("I2").Value = 1282.75
$ws.Range("K2").Value = 1282.75
$ws.Range("M2").Value = -1169.75
$ws.Range("H3").Value = 4764.6665
$ws.Range("I3").Value = 4001.25
$ws.Range("J3").Value = 6291.5
$ws.Range("K3").Value = 4001.25
$ws.Range("L3").Value = 6291.5
$ws.Range("M3").Value = -3886.25
$ws.Range("N3").Value = -6521.5
$ws.Range("H31").Value = 4528.6665
$ws.Range("I31").Value = 4528.6665
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 4528.6665
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -4234.6665
$ws.Range("H32").Value = 27402.928
$ws.Range("I32").Value = 22333.793
$ws.Range("K32").Value = 22333.793
$ws.Range("M32").Value = -22046.793
$ws.Range("H44").Value = 36856.855
$ws.Range("J44").Value = 39666.332
$ws.Range("L44").Value = 39666.332
$ws.Range("N44").Value = -40642.332
$ws.Range("H45").Value = 33335560
$ws.Range("I45").Value = 41668860
$ws.Range("J45").Value = 2361.3333
$ws.Range("K45").Value = 41668860
$ws.Range("L45").Value = 2361.3333
$ws.Range("M45").Value = -41668483
$ws.Range("N45").Value = -3115.3333
$ws.Range("H74").Value = 6677.159
$ws.Range("I74").Value = 2818.6333
$ws.Range("J74").Value = 14945.429
$ws.Range("K74").Value = 2818.6333
$ws.Range("L74").Value = 14945.429
$ws.Range("M74").Value = -1944.6333
$ws.Range("N74").Value = -16693.429
$ws.Range("H77").Value = 6677.159
$ws.Range("I77").Value = 2818.6333
$ws.Range("J77").Value = 14945.429
$ws.Range("K77").Value = 14093.1665
$ws.Range("L77").Value = 74727.145
$ws.Range("M77").Value = -9725.166499999999
$ws.Range("N77").Value = -83463.145
$ws.Range("H116").Value = 1370.4546
$ws.Range("I116").Value = 1282.75
$ws.Range("K116").Value = 1282.75
$ws.Range("M116").Value = 1011.25
$ws.Range("H122").Value = 3023.4827
$ws.Range("I122").Value = 2288.0588
$ws.Range("J122").Value = 4065.3333
$ws.Range("K122").Value = 6864.176399999999
$ws.Range("L122").Value = 12195.9999
$ws.Range("M122").Value = -4414.176399999999
$ws.Range("N122").Value = -17095.9999
$ws.Range("H132").Value = 3299.1904
$ws.Range("I132").Value = 2156.7234
$ws.Range("K132").Value = 6470.1702
$ws.Range("M132").Value = -3940.1702

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1370.4546
$ws.Range("I3").Value = 1282.75
$ws.Range("K3").Value = 1282.75
$ws.Range("M3").Value = -1168.75
$ws.Range("H134").Value = 2269761.2
$ws.Range("I134").Value = 1465.0408
$ws.Range("K134").Value = 4395.1224
$ws.Range("M134").Value = -1860.1224

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 3499.5
$ws.Range("I10").Value = 4928.5
$ws.Range("J10").Value = 641.5
$ws.Range("K10").Value = 4928.5
$ws.Range("L10").Value = 641.5
$ws.Range("M10").Value = -4789.5
$ws.Range("N10").Value = -919.5
$ws.Range("H31").Value = 8665.15
$ws.Range("I31").Value = 9850.904
$ws.Range("J31").Value = 7354.579
$ws.Range("K31").Value = 9850.904
$ws.Range("L31").Value = 7354.579
$ws.Range("M31").Value = -9555.904
$ws.Range("N31").Value = -7944.579
$ws.Range("H34").Value = 8665.15
$ws.Range("I34").Value = 9850.904
$ws.Range("J34").Value = 7354.579
$ws.Range("K34").Value = 9850.904
$ws.Range("L34").Value = 7354.579
$ws.Range("M34").Value = -9648.904
$ws.Range("N34").Value = -7758.579
$ws.Range("H51").Value = 35000
$ws.Range("H61").Value = 35000
$ws.Range("I94").Value = 4069.5833
$ws.Range("J94").Value = 4090.3125
$ws.Range("K94").Value = 4069.5833
$ws.Range("L94").Value = 4090.3125
$ws.Range("M94").Value = -3618.5833
$ws.Range("N94").Value = -4992.3125
$ws.Range("H132").Value = 4040.2559
$ws.Range("I132").Value = 2158.1943
$ws.Range("J132").Value = 13719.429
$ws.Range("K132").Value = 6474.5829
$ws.Range("L132").Value = 41158.287
$ws.Range("M132").Value = -3944.5829
$ws.Range("N132").Value = -46218.287

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 15643.462
$ws.Range("I2").Value = 97.125
$ws.Range("J2").Value = 40517.6
$ws.Range("K2").Value = 582.75
$ws.Range("L2").Value = 243105.6
$ws.Range("M2").Value = -469.75
$ws.Range("N2").Value = -243331.6
$ws.Range("H4").Value = 7747672
$ws.Range("I4").Value = 6500212
$ws.Range("K4").Value = 19500636
$ws.Range("M4").Value = -19500524
$ws.Range("H6").Value = 322.22223
$ws.Range("I6").Value = 30.25
$ws.Range("J6").Value = 555.8
$ws.Range("K6").Value = 90.75
$ws.Range("L6").Value = 1667.4
$ws.Range("M6").Value = 22.25
$ws.Range("N6").Value = -1893.4
$ws.Range("H68").Value = 2748.8572
$ws.Range("I68").Value = 2543
$ws.Range("J68").Value = 2783.1667
$ws.Range("K68").Value = 7629
$ws.Range("L68").Value = 8349.500100000001
$ws.Range("M68").Value = -6818
$ws.Range("N68").Value = -9971.500100000001
$ws.Range("H71").Value = 2748.8572
$ws.Range("I71").Value = 2543
$ws.Range("J71").Value = 2783.1667
$ws.Range("K71").Value = 22887
$ws.Range("L71").Value = 25048.5003
$ws.Range("M71").Value = -18831
$ws.Range("N71").Value = -33160.5003
$ws.Range("H113").Value = 1451.25
$ws.Range("I113").Value = 1040.625
$ws.Range("J113").Value = 1725
$ws.Range("K113").Value = 3121.875
$ws.Range("L113").Value = 5175
$ws.Range("M113").Value = -951.875
$ws.Range("N113").Value = -9515
$ws.Range("H122").Value = 2424.838
$ws.Range("J122").Value = 2611.5151
$ws.Range("L122").Value = 23503.6359
$ws.Range("N122").Value = -28403.6359
$ws.Range("H132").Value = 1677
$ws.Range("J132").Value = 1500
$ws.Range("L132").Value = 13500
$ws.Range("N132").Value = -18560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3038504
$ws.Range("J11").Value = 1258747.2
$ws.Range("L11").Value = 1258747.2
$ws.Range("N11").Value = -1259025.2
$ws.Range("H80").Value = 6736.4
$ws.Range("I80").Value = 3952.5
$ws.Range("K80").Value = 3952.5
$ws.Range("M80").Value = -2954.5
$ws.Range("H83").Value = 6736.4
$ws.Range("I83").Value = 3952.5
$ws.Range("K83").Value = 19762.5
$ws.Range("M83").Value = -14770.5
$ws.Range("H128").Value = 119498.336
$ws.Range("J128").Value = 119498.336
$ws.Range("L128").Value = 119498.336
$ws.Range("N128").Value = -129458.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 19987
$ws.Range("I2").Value = 19950
$ws.Range("K2").Value = 19950
$ws.Range("M2").Value = -19838
$ws.Range("H55").Value = 1137.2858
$ws.Range("I55").Value = 1152.8
$ws.Range("J55").Value = 1098.5
$ws.Range("K55").Value = 1152.8
$ws.Range("L55").Value = 1098.5
$ws.Range("M55").Value = -979.8
$ws.Range("N55").Value = -1444.5
$ws.Range("H68").Value = 4374
$ws.Range("I68").Value = 4061
$ws.Range("K68").Value = 4061
$ws.Range("M68").Value = -3312
$ws.Range("H71").Value = 4374
$ws.Range("I71").Value = 4061
$ws.Range("K71").Value = 20305
$ws.Range("M71").Value = -16561
$ws.Range("H132").Value = 11880.546
$ws.Range("I132").Value = 13659.556
$ws.Range("K132").Value = 40978.66800000001
$ws.Range("M132").Value = -38448.66800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 257500
$ws.Range("I2").Value = 300000
$ws.Range("K2").Value = 300000
$ws.Range("M2").Value = -299888
$ws.Range("H49").Value = 37618
$ws.Range("J49").Value = 38995.668
$ws.Range("L49").Value = 38995.668
$ws.Range("N49").Value = -39455.668
$ws.Range("H122").Value = 3160.6333
$ws.Range("I122").Value = 2963.7036
$ws.Range("K122").Value = 8891.110799999999
$ws.Range("M122").Value = -6441.110799999999
$ws.Range("H132").Value = 3477929
$ws.Range("I132").Value = 6089.0454
$ws.Range("K132").Value = 18267.1362
$ws.Range("M132").Value = -15737.1362
